$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; temporarily unprotect to apply the updates, then
# restore protection afterwards.
$ws.Unprotect()

# Update the confidential notice date from 2021-05-14 to 2021-05-17
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2426975712315811
$ws.Range("E2").Value = -0.007267196214297722

$ws.Range("D3").Value = 0.5030248474321338
$ws.Range("E3").Value = -0.001303780964798107

$ws.Range("D4").Value = 0.09494773737808833
$ws.Range("E4").Value = -0.00899460323805723

$ws.Range("D5").Value = 0.1029142380218218
$ws.Range("E5").Value = -0.001611603545527585

$ws.Range("D6").Value = 0.05641560593637504
$ws.Range("E6").Value = 0.0002267573696144165

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = -0.003426646614316908

$ws.Protect()
